# Updates cryptos list values (Price / Volume(1h) columns) to match
# the latest scrape, and fixes the ordering of two coin-pairs whose
# rows were swapped (PancakeSwap/Binance-PegBSC-USD and Mantle/VeChain).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.496.23"
$ws.Range("E2").Value = "  +1.77%  "

$ws.Range("D3").Value = "'2.653.45"
$ws.Range("E3").Value = "  +1.02%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "'609.92"
$ws.Range("E5").Value = "  +2.56%  "

$ws.Range("D6").Value = "'156.36"
$ws.Range("E6").Value = "  +2.43%  "

$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("D8").Value = "'0.588"
$ws.Range("E8").Value = "  -0.20%  "

$ws.Range("D9").Value = "'2.652.55"
$ws.Range("E9").Value = "  +1.03%  "

$ws.Range("E10").Value = "  +8.38%  "

$ws.Range("D11").Value = "'0.403"
$ws.Range("E11").Value = "  +2.02%  "

$ws.Range("D12").Value = "'5.89"
$ws.Range("E12").Value = "  +1.75%  "

$ws.Range("E13").Value = "  +1.67%  "

$ws.Range("E14").Value = "  +5.54%  "

$ws.Range("E15").Value = "  +15.41%  "

$ws.Range("D16").Value = "'3.131.00"
$ws.Range("E16").Value = "  +1.04%  "

$ws.Range("D17").Value = "'65.283.21"
$ws.Range("E17").Value = "  +1.54%  "

$ws.Range("D18").Value = "'2.657.68"
$ws.Range("E18").Value = "  +0.70%  "

$ws.Range("D19").Value = "'12.75"
$ws.Range("E19").Value = "  +3.91%  "

$ws.Range("D20").Value = "'4.90"
$ws.Range("E20").Value = "  +2.88%  "

$ws.Range("D21").Value = "'359.20"
$ws.Range("E21").Value = "  +3.07%  "

$ws.Range("D22").Value = "'7.44"
$ws.Range("E22").Value = "  +4.52%  "

$ws.Range("E23").Value = "  +0.16%  "

$ws.Range("D24").Value = "'69.90"
$ws.Range("E24").Value = "  +3.40%  "

$ws.Range("D25").Value = "'1.72"
$ws.Range("E25").Value = "  +1.82%  "

$ws.Range("D26").Value = "'9.48"
$ws.Range("E26").Value = "  +2.57%  "

$ws.Range("D27").Value = "'0.0000106"
$ws.Range("E27").Value = "  +16.69%  "

$ws.Range("D28").Value = "'1.64"
$ws.Range("E28").Value = "  -1.37%  "

$ws.Range("E29").Value = "  +2.42%  "

$ws.Range("D30").Value = "'8.17"
$ws.Range("E30").Value = "  -0.56%  "

$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'2.19"
$ws.Range("E31").Value = "  +5.34%  "

$ws.Range("B32").Value = "Binance-PegBSC-USD"
$ws.Range("C32").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D32").Value = "'1.00"
$ws.Range("E32").Value = "  -0.06%  "

$ws.Range("D33").Value = "'526.13"
$ws.Range("E33").Value = "  -4.14%  "

$ws.Range("D34").Value = "'1.79"
$ws.Range("E34").Value = "  -1.00%  "

$ws.Range("D35").Value = "'5.54"
$ws.Range("E35").Value = "  +0.19%  "

$ws.Range("D36").Value = "'6.37"
$ws.Range("E36").Value = "  +2.98%  "

$ws.Range("D37").Value = "'0.432"
$ws.Range("E37").Value = "  +2.83%  "

$ws.Range("D38").Value = "'20.71"
$ws.Range("E38").Value = "  +3.53%  "

$ws.Range("D39").Value = "'162.83"
$ws.Range("E39").Value = "  -1.08%  "

$ws.Range("D40").Value = "'1.98"
$ws.Range("E40").Value = "  -0.73%  "

$ws.Range("E42").Value = "  +0.01%  "

$ws.Range("D43").Value = "'41.95"
$ws.Range("E43").Value = "  +0.78%  "

$ws.Range("D44").Value = "'165.76"
$ws.Range("E44").Value = "  -1.16%  "

$ws.Range("D45").Value = "'4.14"
$ws.Range("E45").Value = "  +0.77%  "

$ws.Range("D46").Value = "'2.35"
$ws.Range("E46").Value = "  +5.91%  "

$ws.Range("D47").Value = "'0.0613"
$ws.Range("E47").Value = "  +4.46%  "

$ws.Range("D48").Value = "'22.98"
$ws.Range("E48").Value = "  -1.89%  "

$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").Value = "'0.0264"
$ws.Range("E49").Value = "  +5.30%  "

$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").Value = "'0.653"
$ws.Range("E50").Value = "  +1.76%  "

$ws.Range("D51").Value = "'0.0981"
$ws.Range("E51").Value = "  +0.29%  "
